# Apply the commit's edits:
#  - Sheet2!A2: formula changed from =Sheet1!B12/1000000 to =Sheet1!B11*100
#    (and its number format reset to General, which is what produced the
#    extra cellXfs entry in the saved styles part)
#  - Sheet1 selection moved to B11
#  - Sheet2 selection moved to A3

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update the LVaR formula on Sheet2
$ws2.Range("A2").Formula = "=Sheet1!B11*100"
$ws2.Range("A2").NumberFormat = "General"

# Restore the active selections on each sheet
$ws1.Range("B11").Select()
$ws2.Range("A3").Select()
